# Natmi following Dr Hou advice
# Updates LR-pair statistics for Npnt-Itga8 (rows 2-13) and appends new
# rows 14-16 to extend the Target cluster coverage from 5 to 6 groups per
# Sending cluster (ECs/FAPs/sCs), each against Npnt/Itga8 x {ECs,FAPs,M1,M2,sCs}.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Npnt"
$ws.Cells.Item(2, 3).Value = "Itga8"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.440918333333333
$ws.Cells.Item(2, 8).Value = 7.322755
$ws.Cells.Item(2, 9).Value = 0.5182826554654038
$ws.Cells.Item(2, 10).Value = 0.5182826554654038
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.3881176666666666
$ws.Cells.Item(2, 14).Value = 1.164353
$ws.Cells.Item(2, 15).Value = 0.08027910179751364
$ws.Cells.Item(2, 16).Value = 0.08027910179751364
$ws.Cells.Item(2, 17).Value = 0.9473635280572221
$ws.Cells.Item(2, 18).Value = 8.526271752515
$ws.Cells.Item(2, 19).Value = 0.04160726605799284
$ws.Cells.Item(2, 20).Value = 0.04160726605799284
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Npnt"
$ws.Cells.Item(3, 3).Value = "Itga8"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.440918333333333
$ws.Cells.Item(3, 8).Value = 7.322755
$ws.Cells.Item(3, 9).Value = 0.5182826554654038
$ws.Cells.Item(3, 10).Value = 0.5182826554654038
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.230730666666667
$ws.Cells.Item(3, 14).Value = 6.692192
$ws.Cells.Item(3, 15).Value = 0.4614091798762974
$ws.Cells.Item(3, 16).Value = 0.4614091798762974
$ws.Cells.Item(3, 17).Value = 5.445031380995555
$ws.Cells.Item(3, 18).Value = 49.00528242896
$ws.Cells.Item(3, 19).Value = 0.2391403750024015
$ws.Cells.Item(3, 20).Value = 0.2391403750024015
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Npnt"
$ws.Cells.Item(4, 3).Value = "Itga8"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.440918333333333
$ws.Cells.Item(4, 8).Value = 7.322755
$ws.Cells.Item(4, 9).Value = 0.5182826554654038
$ws.Cells.Item(4, 10).Value = 0.5182826554654038
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.9049766666666667
$ws.Cells.Item(4, 14).Value = 2.71493
$ws.Cells.Item(4, 15).Value = 0.1871873408177105
$ws.Cells.Item(4, 16).Value = 0.1871873408177105
$ws.Cells.Item(4, 17).Value = 2.208974136905555
$ws.Cells.Item(4, 18).Value = 19.88076723215
$ws.Cells.Item(4, 19).Value = 0.09701595206851059
$ws.Cells.Item(4, 20).Value = 0.09701595206851059
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Npnt"
$ws.Cells.Item(5, 3).Value = "Itga8"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.440918333333333
$ws.Cells.Item(5, 8).Value = 7.322755
$ws.Cells.Item(5, 9).Value = 0.5182826554654038
$ws.Cells.Item(5, 10).Value = 0.5182826554654038
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.6326649999999999
$ws.Cells.Item(5, 14).Value = 1.897995
$ws.Cells.Item(5, 15).Value = 0.1308618037795857
$ws.Cells.Item(5, 16).Value = 0.1308618037795857
$ws.Cells.Item(5, 17).Value = 1.544283597358333
$ws.Cells.Item(5, 18).Value = 13.898552376225
$ws.Cells.Item(5, 19).Value = 0.06782340316187627
$ws.Cells.Item(5, 20).Value = 0.06782340316187627
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Npnt"
$ws.Cells.Item(6, 3).Value = "Itga8"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.440918333333333
$ws.Cells.Item(6, 8).Value = 7.322755
$ws.Cells.Item(6, 9).Value = 0.5182826554654038
$ws.Cells.Item(6, 10).Value = 0.5182826554654038
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.678114
$ws.Cells.Item(6, 14).Value = 2.034342
$ws.Cells.Item(6, 15).Value = 0.1402625737288928
$ws.Cells.Item(6, 16).Value = 0.1402625737288928
$ws.Cells.Item(6, 17).Value = 1.65522089469
$ws.Cells.Item(6, 18).Value = 14.89698805221
$ws.Cells.Item(6, 19).Value = 0.07269565917462255
$ws.Cells.Item(6, 20).Value = 0.07269565917462255
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Npnt"
$ws.Cells.Item(7, 3).Value = "Itga8"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.6061233333333333
$ws.Cells.Item(7, 8).Value = 1.81837
$ws.Cells.Item(7, 9).Value = 0.1286987796558298
$ws.Cells.Item(7, 10).Value = 0.1286987796558298
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.3881176666666666
$ws.Cells.Item(7, 14).Value = 1.164353
$ws.Cells.Item(7, 15).Value = 0.08027910179751364
$ws.Cells.Item(7, 16).Value = 0.08027910179751364
$ws.Cells.Item(7, 17).Value = 0.2352471738455555
$ws.Cells.Item(7, 18).Value = 2.11722456461
$ws.Cells.Item(7, 19).Value = 0.01033182243320614
$ws.Cells.Item(7, 20).Value = 0.01033182243320614
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Npnt"
$ws.Cells.Item(8, 3).Value = "Itga8"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.6061233333333333
$ws.Cells.Item(8, 8).Value = 1.81837
$ws.Cells.Item(8, 9).Value = 0.1286987796558298
$ws.Cells.Item(8, 10).Value = 0.1286987796558298
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.230730666666667
$ws.Cells.Item(8, 14).Value = 6.692192
$ws.Cells.Item(8, 15).Value = 0.4614091798762974
$ws.Cells.Item(8, 16).Value = 0.4614091798762974
$ws.Cells.Item(8, 17).Value = 1.352097907448889
$ws.Cells.Item(8, 18).Value = 12.16888116704
$ws.Cells.Item(8, 19).Value = 0.05938279837207674
$ws.Cells.Item(8, 20).Value = 0.05938279837207675
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Npnt"
$ws.Cells.Item(9, 3).Value = "Itga8"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.6061233333333333
$ws.Cells.Item(9, 8).Value = 1.81837
$ws.Cells.Item(9, 9).Value = 0.1286987796558298
$ws.Cells.Item(9, 10).Value = 0.1286987796558298
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.9049766666666667
$ws.Cells.Item(9, 14).Value = 2.71493
$ws.Cells.Item(9, 15).Value = 0.1871873408177105
$ws.Cells.Item(9, 16).Value = 0.1871873408177105
$ws.Cells.Item(9, 17).Value = 0.5485274737888889
$ws.Cells.Item(9, 18).Value = 4.9367472641
$ws.Cells.Item(9, 19).Value = 0.02409078233025925
$ws.Cells.Item(9, 20).Value = 0.02409078233025926
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Npnt"
$ws.Cells.Item(10, 3).Value = "Itga8"
$ws.Cells.Item(10, 4).Value = "M2"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.6061233333333333
$ws.Cells.Item(10, 8).Value = 1.81837
$ws.Cells.Item(10, 9).Value = 0.1286987796558298
$ws.Cells.Item(10, 10).Value = 0.1286987796558298
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.6326649999999999
$ws.Cells.Item(10, 14).Value = 1.897995
$ws.Cells.Item(10, 15).Value = 0.1308618037795857
$ws.Cells.Item(10, 16).Value = 0.1308618037795857
$ws.Cells.Item(10, 17).Value = 0.3834730186833333
$ws.Cells.Item(10, 18).Value = 3.45125716815
$ws.Cells.Item(10, 19).Value = 0.01684175444999333
$ws.Cells.Item(10, 20).Value = 0.01684175444999334
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Npnt"
$ws.Cells.Item(11, 3).Value = "Itga8"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.6061233333333333
$ws.Cells.Item(11, 8).Value = 1.81837
$ws.Cells.Item(11, 9).Value = 0.1286987796558298
$ws.Cells.Item(11, 10).Value = 0.1286987796558298
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.678114
$ws.Cells.Item(11, 14).Value = 2.034342
$ws.Cells.Item(11, 15).Value = 0.1402625737288928
$ws.Cells.Item(11, 16).Value = 0.1402625737288928
$ws.Cells.Item(11, 17).Value = 0.41102071806
$ws.Cells.Item(11, 18).Value = 3.69918646254
$ws.Cells.Item(11, 19).Value = 0.01805162207029436
$ws.Cells.Item(11, 20).Value = 0.01805162207029436
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Npnt"
$ws.Cells.Item(12, 3).Value = "Itga8"
$ws.Cells.Item(12, 4).Value = "ECs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.662586
$ws.Cells.Item(12, 8).Value = 4.987757999999999
$ws.Cells.Item(12, 9).Value = 0.3530185648787664
$ws.Cells.Item(12, 10).Value = 0.3530185648787664
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.3881176666666666
$ws.Cells.Item(12, 14).Value = 1.164353
$ws.Cells.Item(12, 15).Value = 0.08027910179751364
$ws.Cells.Item(12, 16).Value = 0.08027910179751364
$ws.Cells.Item(12, 17).Value = 0.6452789989526666
$ws.Cells.Item(12, 18).Value = 5.807510990573999
$ws.Cells.Item(12, 19).Value = 0.02834001330631466
$ws.Cells.Item(12, 20).Value = 0.02834001330631466
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Npnt"
$ws.Cells.Item(13, 3).Value = "Itga8"
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.662586
$ws.Cells.Item(13, 8).Value = 4.987757999999999
$ws.Cells.Item(13, 9).Value = 0.3530185648787664
$ws.Cells.Item(13, 10).Value = 0.3530185648787664
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.230730666666667
$ws.Cells.Item(13, 14).Value = 6.692192
$ws.Cells.Item(13, 15).Value = 0.4614091798762974
$ws.Cells.Item(13, 16).Value = 0.4614091798762974
$ws.Cells.Item(13, 17).Value = 3.708781576170666
$ws.Cells.Item(13, 18).Value = 33.379034185536
$ws.Cells.Item(13, 19).Value = 0.1628860065018191
$ws.Cells.Item(13, 20).Value = 0.1628860065018191
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Npnt"
$ws.Cells.Item(14, 3).Value = "Itga8"
$ws.Cells.Item(14, 4).Value = "M1"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1.662586
$ws.Cells.Item(14, 8).Value = 4.987757999999999
$ws.Cells.Item(14, 9).Value = 0.3530185648787664
$ws.Cells.Item(14, 10).Value = 0.3530185648787664
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.9049766666666667
$ws.Cells.Item(14, 14).Value = 2.71493
$ws.Cells.Item(14, 15).Value = 0.1871873408177105
$ws.Cells.Item(14, 16).Value = 0.1871873408177105
$ws.Cells.Item(14, 17).Value = 1.504601536326666
$ws.Cells.Item(14, 18).Value = 13.54141382694
$ws.Cells.Item(14, 19).Value = 0.0660806064189407
$ws.Cells.Item(14, 20).Value = 0.0660806064189407
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Npnt"
$ws.Cells.Item(15, 3).Value = "Itga8"
$ws.Cells.Item(15, 4).Value = "M2"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 1.662586
$ws.Cells.Item(15, 8).Value = 4.987757999999999
$ws.Cells.Item(15, 9).Value = 0.3530185648787664
$ws.Cells.Item(15, 10).Value = 0.3530185648787664
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.6326649999999999
$ws.Cells.Item(15, 14).Value = 1.897995
$ws.Cells.Item(15, 15).Value = 0.1308618037795857
$ws.Cells.Item(15, 16).Value = 0.1308618037795857
$ws.Cells.Item(15, 17).Value = 1.05185997169
$ws.Cells.Item(15, 18).Value = 9.466739745209999
$ws.Cells.Item(15, 19).Value = 0.04619664616771606
$ws.Cells.Item(15, 20).Value = 0.04619664616771606
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Npnt"
$ws.Cells.Item(16, 3).Value = "Itga8"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 1.662586
$ws.Cells.Item(16, 8).Value = 4.987757999999999
$ws.Cells.Item(16, 9).Value = 0.3530185648787664
$ws.Cells.Item(16, 10).Value = 0.3530185648787664
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.678114
$ws.Cells.Item(16, 14).Value = 2.034342
$ws.Cells.Item(16, 15).Value = 0.1402625737288928
$ws.Cells.Item(16, 16).Value = 0.1402625737288928
$ws.Cells.Item(16, 17).Value = 1.127422842804
$ws.Cells.Item(16, 18).Value = 10.146805585236
$ws.Cells.Item(16, 19).Value = 0.0495152924839759
$ws.Cells.Item(16, 20).Value = 0.0495152924839759
